$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B6").Value = 'Liste 1 et 2'
$ws.Range("C6").Value = 'enlever liste 1 et 2 sur annuaires'
$ws.Range("D6").Value = 'ça ne sert à rien'
$ws.Range("E6").Value = 'retirer liste 1 et 2'

$ws.Range("B10").Value = 'le point sur H1'
$ws.Range("C10").Value = 'enlever le point à la fin du titre'
$ws.Range("D10").Value = 'il n’y a pas besoin de point à la fin du titre'
$ws.Range("E10").Value = 'enlever le point '

$ws.Range("B13").Value = 'images'
$ws.Range("C13").Value = 'changer le nom des images en jpeg'
$ws.Range("D13").Value = 'cela prend trop de temps pour le navigateur de les télécharger'
$ws.Range("E13").Value = 'renommer les images'

$ws.Range("B14").Value = 'et-line.min'
$ws.Range("C14").Value = 'retirer .min '
$ws.Range("D14").Value = 'le et-line.min n’existe pas les modifications ne s’appliquent pas'
$ws.Range("E14").Value = 'retirer .min'

$ws.Range("B15").Value = 'Bootstrap.min sur contact.html'
$ws.Range("C15").Value = 'retirer .min '
$ws.Range("D15").Value = 'le bootstrap.min n’existe pas les modifications ne s’appliquent pas'
$ws.Range("E15").Value = 'retirer .min'

$ws.Range("B16").Value = 'Font-awesome.min sur contact.html'
$ws.Range("C16").Value = 'retirer .min '
$ws.Range("D16").Value = 'le font-awesome n’existe pas les modifications ne s’appliquent pas'
$ws.Range("E16").Value = 'retirer .min'

$ws.Range("B17").Value = 'couleur liste annuaires et partenaires'
$ws.Range("C17").Value = 'changer la couleur des liens sur les listes'
$ws.Range("D17").Value = 'on ne voit pas la différence entre les titres et les liens'
$ws.Range("E17").Value = 'changer la couleur '

$ws.Columns.Item(4).ColumnWidth = 54.38
$ws.Range("E17").Select() | Out-Null
